$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D holds numeric-looking values stored as TEXT (inline strings) in
#     the workbook.  If we just assign a numeric-looking string, Excel will
#     coerce it into a real number (losing the exact textual representation
#     and changing the cell type).  To keep these as text we temporarily mark
#     the cell as Text ("@") before assigning the value, then restore the
#     cell's style afterwards so no visual/style changes leak into the file.

$priceUpdates = @{
  "D2"  = "246.74"
  "D3"  = "22.43"
  "D4"  = "5.485"
  "D5"  = "0.05638"
  "D6"  = "6.466"
  "D7"  = "0.8064"
  "D8"  = "1.045"
  "D9"  = "0.1444"
  "D10" = "0.07368"
  "D11" = "0.03194"
  "D12" = "0.02933"
  "D13" = "0.09257"
  "D14" = "0.001670"
  "D15" = "3.202"
  "D16" = "0.04721"
  "D17" = "0.0005836"
  "D18" = "0.006313"
  "D19" = "0.001053"
  "D20" = "0.004115"
  "D21" = "0.0001505"
  "D22" = "3.980"
  "D23" = "3.384"
  "D26" = "0.1314"
  "D27" = "0.0003008"
  "D40" = "0.04154"
  "D41" = "0.006883"
  "D42" = "0.003511"
  "D43" = "0.1040"
  "D44" = "0.009051"
  "D45" = "0.00005645"
  "D46" = "0.00000000752"
  "D47" = "0.6819"
  "D48" = "0.01995"
  "D49" = "0.00002107"
  "D50" = "0.01013"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# --- Row 41 and 43 coin swap (BKEXToken <-> KickToken), plus their
#     accompanying link / 24h-volume-label text updates.

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Row 17 "Worst in 24h" label addition.
$ws.Range("E17").Value = "16OneONEWorstin24h"
